# Fruta / hortaliza, semanal
# Insert two new weekly observation rows into the Pimiento sheet just
# before the existing row 414, pushing the former rows 414-434 down to
# 416-436, then populate the two freshly-inserted rows (414-415) with
# their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 414 (old rows 414-434 shift to 416-436).
$ws.Rows("414:415").Insert()

# New row 414: Zafiro rojo
$ws.Cells.Item(414, 1).Value2 = 7
$ws.Cells.Item(414, 2).Value2 = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(414, 3).Value2 = 'Ñuble'
$ws.Cells.Item(414, 4).Value2 = 45013
$ws.Cells.Item(414, 5).Value2 = 16
$ws.Cells.Item(414, 6).Value2 = 100112002
$ws.Cells.Item(414, 7).Value2 = 'Pimiento'
$ws.Cells.Item(414, 8).Value2 = 'Zafiro rojo'
$ws.Cells.Item(414, 9).Value2 = 'Primera'
$ws.Cells.Item(414, 10).Value2 = 40
$ws.Cells.Item(414, 11).Value2 = 17000
$ws.Cells.Item(414, 12).Value2 = 17000
$ws.Cells.Item(414, 13).Value2 = 17000
$ws.Cells.Item(414, 14).Value2 = '$/caja 15 kilos'
$ws.Cells.Item(414, 15).Value2 = 'Región de Arica y Parinacota'
$ws.Cells.Item(414, 16).Value2 = 1133
$ws.Cells.Item(414, 17).Value2 = 15
$ws.Cells.Item(414, 18).Value2 = 'Hortaliza'

# New row 415: Zafiro verde
$ws.Cells.Item(415, 1).Value2 = 7
$ws.Cells.Item(415, 2).Value2 = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(415, 3).Value2 = 'Ñuble'
$ws.Cells.Item(415, 4).Value2 = 45013
$ws.Cells.Item(415, 5).Value2 = 16
$ws.Cells.Item(415, 6).Value2 = 100112002
$ws.Cells.Item(415, 7).Value2 = 'Pimiento'
$ws.Cells.Item(415, 8).Value2 = 'Zafiro verde'
$ws.Cells.Item(415, 9).Value2 = 'Primera'
$ws.Cells.Item(415, 10).Value2 = 40
$ws.Cells.Item(415, 11).Value2 = 13000
$ws.Cells.Item(415, 12).Value2 = 13000
$ws.Cells.Item(415, 13).Value2 = 13000
$ws.Cells.Item(415, 14).Value2 = '$/caja 15 kilos'
$ws.Cells.Item(415, 15).Value2 = 'Región de Arica y Parinacota'
$ws.Cells.Item(415, 16).Value2 = 867
$ws.Cells.Item(415, 17).Value2 = 15
$ws.Cells.Item(415, 18).Value2 = 'Hortaliza'
